$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 updates
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 5.9
$ws.Range("K6").Value = 2.32
$ws.Range("L6").Value = 5.6
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 3.98
$ws.Range("Q6").Value = 1.6
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.92
$ws.Range("W6").Value = 6.6
$ws.Range("X6").Value = 6.3
$ws.Range("Y6").Value = 6.9
$ws.Range("Z6").Value = 8.5
$ws.Range("AA6").Value = 9.25
$ws.Range("AB6").Value = 18.5
$ws.Range("AD6").Value = 7.2
$ws.Range("AE6").Value = 14
$ws.Range("AH6").Value = 14.5
$ws.Range("AI6").Value = 30
$ws.Range("AJ6").Value = 15.5
$ws.Range("AK6").Value = 90
$ws.Range("AN6").Value = 3.35
$ws.Range("AS6").Value = 175
$ws.Range("AT6").Value = 3.05
$ws.Range("AU6").Value = 7.7
$ws.Range("AV6").Value = 65
$ws.Range("AW6").Value = 7.5
$ws.Range("AX6").Value = 35
$ws.Range("BA6").Value = 250

# Row 9 updates
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.67
$ws.Range("R9").Value = 2.15
